$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("N6").Value = 6.65
$ws.Range("N7").Value = 7.8
